$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows for CRB13, CRB62, CRB72 (delete highest row index first to avoid renumbering issues)
$ws.Rows(29).Delete()
$ws.Rows(26).Delete()
$ws.Rows(4).Delete()

# Update the remaining data grid (columns B:J) with the new values
# Row 2: CRB11
$ws.Range("B2").Value = 4
$ws.Range("C2").Value = 5
$ws.Range("D2").Value = 3
$ws.Range("E2").Value = 11
$ws.Range("F2").Value = 13
$ws.Range("G2").Value = 8
$ws.Range("H2").Value = 10
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 2

# Row 3: CRB12
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 2
$ws.Range("G3").Value = 2
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0

# Row 4: CRB15
$ws.Range("B4").Value = 4
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 3
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 6
$ws.Range("H4").Value = 2
$ws.Range("I4").Value = 2
$ws.Range("J4").Value = 3

# Row 5: CRB21
$ws.Range("B5").Value = 17
$ws.Range("C5").Value = 13
$ws.Range("D5").Value = 16
$ws.Range("E5").Value = 18
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 18
$ws.Range("H5").Value = 16
$ws.Range("I5").Value = 15
$ws.Range("J5").Value = 16

# Row 6: CRB22
$ws.Range("B6").Value = 3
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 3
$ws.Range("E6").Value = 5
$ws.Range("F6").Value = 6
$ws.Range("G6").Value = 1
$ws.Range("H6").Value = 4
$ws.Range("I6").Value = 2
$ws.Range("J6").Value = 5

# Row 7: CRB23
$ws.Range("B7").Value = 2
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 2
$ws.Range("G7").Value = 5
$ws.Range("H7").Value = 6
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 2

# Row 8: CRB24
$ws.Range("B8").Value = 11
$ws.Range("C8").Value = 6
$ws.Range("D8").Value = 8
$ws.Range("E8").Value = 11
$ws.Range("F8").Value = 4
$ws.Range("G8").Value = 10
$ws.Range("H8").Value = 12
$ws.Range("I8").Value = 5
$ws.Range("J8").Value = 9

# Row 9: CRB25
$ws.Range("B9").Value = 4
$ws.Range("C9").Value = 3
$ws.Range("D9").Value = 7
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 5
$ws.Range("G9").Value = 11
$ws.Range("H9").Value = 10
$ws.Range("I9").Value = 14
$ws.Range("J9").Value = 9

# Row 10: CRB26
$ws.Range("B10").Value = 7
$ws.Range("C10").Value = 4
$ws.Range("D10").Value = 5
$ws.Range("E10").Value = 12
$ws.Range("F10").Value = 6
$ws.Range("G10").Value = 15
$ws.Range("H10").Value = 19
$ws.Range("I10").Value = 1
$ws.Range("J10").Value = 7

# Row 11: CRB31
$ws.Range("B11").Value = 11
$ws.Range("C11").Value = 8
$ws.Range("D11").Value = 8
$ws.Range("E11").Value = 7
$ws.Range("F11").Value = 12
$ws.Range("G11").Value = 9
$ws.Range("H11").Value = 6
$ws.Range("I11").Value = 1
$ws.Range("J11").Value = 3

# Row 12: CRB32
$ws.Range("B12").Value = 11
$ws.Range("C12").Value = 3
$ws.Range("D12").Value = 10
$ws.Range("E12").Value = 5
$ws.Range("F12").Value = 8
$ws.Range("G12").Value = 9
$ws.Range("H12").Value = 9
$ws.Range("I12").Value = 6
$ws.Range("J12").Value = 6

# Row 13: CRB34
$ws.Range("B13").Value = 4
$ws.Range("C13").Value = 2
$ws.Range("D13").Value = 4
$ws.Range("E13").Value = 4
$ws.Range("F13").Value = 7
$ws.Range("G13").Value = 8
$ws.Range("H13").Value = 4
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 2

# Row 14: CRB35
$ws.Range("B14").Value = 0
$ws.Range("C14").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 2
$ws.Range("G14").Value = 2
$ws.Range("H14").Value = 1
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 1

# Row 15: CRB41
$ws.Range("B15").Value = 10
$ws.Range("C15").Value = 5
$ws.Range("D15").Value = 11
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 7
$ws.Range("G15").Value = 6
$ws.Range("H15").Value = 13
$ws.Range("I15").Value = 1
$ws.Range("J15").Value = 7

# Row 16: CRB42
$ws.Range("B16").Value = 2
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 2
$ws.Range("G16").Value = 3
$ws.Range("H16").Value = 2
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 3

# Row 17: CRB43
$ws.Range("B17").Value = 10
$ws.Range("C17").Value = 7
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 4
$ws.Range("F17").Value = 5
$ws.Range("G17").Value = 11
$ws.Range("H17").Value = 7
$ws.Range("I17").Value = 2
$ws.Range("J17").Value = 6

# Row 18: CRB44
$ws.Range("B18").Value = 17
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 10
$ws.Range("E18").Value = 6
$ws.Range("F18").Value = 3
$ws.Range("G18").Value = 18
$ws.Range("H18").Value = 17
$ws.Range("I18").Value = 4
$ws.Range("J18").Value = 7

# Row 19: CRB51
$ws.Range("B19").Value = 10
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 18
$ws.Range("E19").Value = 1
$ws.Range("F19").Value = 15
$ws.Range("G19").Value = 19
$ws.Range("H19").Value = 18
$ws.Range("I19").Value = 18
$ws.Range("J19").Value = 19

# Row 20: CRB52
$ws.Range("B20").Value = 2
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = 1
$ws.Range("F20").Value = 3
$ws.Range("G20").Value = 2
$ws.Range("H20").Value = 2
$ws.Range("I20").Value = 1
$ws.Range("J20").Value = 3

# Row 21: CRB53
$ws.Range("B21").Value = 1
$ws.Range("C21").Value = 1
$ws.Range("D21").Value = 0
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = 2
$ws.Range("H21").Value = 1
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 0

# Row 22: CRB54
$ws.Range("B22").Value = 3
$ws.Range("C22").Value = 1
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 1
$ws.Range("F22").Value = 5
$ws.Range("G22").Value = 6
$ws.Range("H22").Value = 16
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 2

# Row 23: CRB55
$ws.Range("B23").Value = 0
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 1
$ws.Range("F23").Value = 3
$ws.Range("G23").Value = 5
$ws.Range("H23").Value = 1
$ws.Range("I23").Value = 4
$ws.Range("J23").Value = 0

# Row 24: CRB61
$ws.Range("B24").Value = 11
$ws.Range("C24").Value = 10
$ws.Range("D24").Value = 11
$ws.Range("E24").Value = 11
$ws.Range("F24").Value = 11
$ws.Range("G24").Value = 12
$ws.Range("H24").Value = 11
$ws.Range("I24").Value = 4
$ws.Range("J24").Value = 0

# Row 25: CRB65
$ws.Range("B25").Value = 7
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 2
$ws.Range("E25").Value = 5
$ws.Range("F25").Value = 5
$ws.Range("G25").Value = 8
$ws.Range("H25").Value = 8
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 1

# Row 26: CRB71
$ws.Range("B26").Value = 6
$ws.Range("C26").Value = 8
$ws.Range("D26").Value = 5
$ws.Range("E26").Value = 7
$ws.Range("F26").Value = 2
$ws.Range("G26").Value = 17
$ws.Range("H26").Value = 16
$ws.Range("I26").Value = 5
$ws.Range("J26").Value = 2

# Row 27: CRB73
$ws.Range("B27").Value = 3
$ws.Range("C27").Value = 0
$ws.Range("D27").Value = 0
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 0
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = 2
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0

# Row 28: CRB74
$ws.Range("B28").Value = 5
$ws.Range("C28").Value = 2
$ws.Range("D28").Value = 7
$ws.Range("E28").Value = 2
$ws.Range("F28").Value = 1
$ws.Range("G28").Value = 13
$ws.Range("H28").Value = 18
$ws.Range("I28").Value = 2
$ws.Range("J28").Value = 7

# Row 29: CRB75
$ws.Range("B29").Value = 14
$ws.Range("C29").Value = 12
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 14
$ws.Range("F29").Value = 1
$ws.Range("G29").Value = 9
$ws.Range("H29").Value = 12
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 2
